$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 5
$ws.Range("D2").Value = 0

$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 0
